$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 17 (ENW017 test case): prepend the new OPQA-3575 jira id to the existing
# B17 ids, and replace the description in C17 with the combined text that now
# also covers the "community enabled version of Endnote" verification step.
$ws.Range("B17").Value = "OPQA-3575`n||OPQA-2157 `n||OPQA-2159 "
$ws.Range("C17").Value = "Verify that user shall be sent to Community enabled version of Endnote while user sign in to the ENW through STeAM or Social as a first time,when the user is affiliated to a Customer in the market test group based on the WOS Customer Check.`n||Verify that the  Endnote profile fly out should contain `"Feedback`" link, which should take the user to the Neon version of the Endnote Feedback form.               `n || And Verify that,the user's message should be sent to a configurable email box specific for Endnote, when user submitting a message in the help input form on the Endnote version of the new `"Feedback`" page"

# The extra line added to the description makes row 17 taller.
$ws.Rows(17).RowHeight = 120

# Leave the selection on the cell that was edited.
$ws.Range("C17").Select()
